# Update NATMI ligand-receptor pair TPM-derived metrics (F7-F3) with
# freshly recomputed values (new TPM input), per "update scripts wuth new tpm".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.012198
$ws.Range("H2").Value = 9.036594000000001
$ws.Range("I2").Value = 0.7974910863989846
$ws.Range("J2").Value = 0.7974910863989846
$ws.Range("M2").Value = 0.2054816666666667
$ws.Range("N2").Value = 0.616445
$ws.Range("O2").Value = 0.00496424614546655
$ws.Range("P2").Value = 0.004964246145466549
$ws.Range("Q2").Value = 0.61895146537
$ws.Range("R2").Value = 5.57056318833
$ws.Range("S2").Value = 0.00395894205170009
$ws.Range("T2").Value = 0.00395894205170009
$ws.Range("G3").Value = 3.012198
$ws.Range("H3").Value = 9.036594000000001
$ws.Range("I3").Value = 0.7974910863989846
$ws.Range("J3").Value = 0.7974910863989846
$ws.Range("O3").Value = 0.9529850468799925
$ws.Range("P3").Value = 0.9529850468799924
$ws.Range("Q3").Value = 118.819952507942
$ws.Range("R3").Value = 1069.379572571478
$ws.Range("S3").Value = 0.7599970803583125
$ws.Range("T3").Value = 0.7599970803583124
$ws.Range("G4").Value = 3.012198
$ws.Range("H4").Value = 9.036594000000001
$ws.Range("I4").Value = 0.7974910863989846
$ws.Range("J4").Value = 0.7974910863989846
$ws.Range("M4").Value = 1.712817
$ws.Range("N4").Value = 5.138451
$ws.Range("O4").Value = 0.04138006727350978
$ws.Range("P4").Value = 0.04138006727350978
$ws.Range("Q4").Value = 5.159343941766
$ws.Range("R4").Value = 46.434095475894
$ws.Range("S4").Value = 0.03300023480521438
$ws.Range("T4").Value = 0.03300023480521438
$ws.Range("G5").Value = 3.012198
$ws.Range("H5").Value = 9.036594000000001
$ws.Range("I5").Value = 0.7974910863989846
$ws.Range("J5").Value = 0.7974910863989846
$ws.Range("M5").Value = 0.02775933333333333
$ws.Range("N5").Value = 0.083278
$ws.Range("O5").Value = 0.000670639701031176
$ws.Range("P5").Value = 0.000670639701031176
$ws.Range("Q5").Value = 0.083616608348
$ws.Range("R5").Value = 0.7525494751320001
$ws.Range("S5").Value = 0.0005348291837576428
$ws.Range("T5").Value = 0.0005348291837576428
$ws.Range("G6").Value = 0.7648950000000001
$ws.Range("H6").Value = 2.294685
$ws.Range("I6").Value = 0.2025089136010154
$ws.Range("J6").Value = 0.2025089136010154
$ws.Range("M6").Value = 0.2054816666666667
$ws.Range("N6").Value = 0.616445
$ws.Range("O6").Value = 0.00496424614546655
$ws.Range("P6").Value = 0.004964246145466549
$ws.Range("Q6").Value = 0.157171899425
$ws.Range("R6").Value = 1.414547094825
$ws.Range("S6").Value = 0.001005304093766459
$ws.Range("T6").Value = 0.001005304093766459
$ws.Range("G7").Value = 0.7648950000000001
$ws.Range("H7").Value = 2.294685
$ws.Range("I7").Value = 0.2025089136010154
$ws.Range("J7").Value = 0.2025089136010154
$ws.Range("O7").Value = 0.9529850468799925
$ws.Range("P7").Value = 0.9529850468799924
$ws.Range("Q7").Value = 30.172248827455
$ws.Range("R7").Value = 271.550239447095
$ws.Range("S7").Value = 0.19298796652168
$ws.Range("T7").Value = 0.19298796652168
$ws.Range("G8").Value = 0.7648950000000001
$ws.Range("H8").Value = 2.294685
$ws.Range("I8").Value = 0.2025089136010154
$ws.Range("J8").Value = 0.2025089136010154
$ws.Range("M8").Value = 1.712817
$ws.Range("N8").Value = 5.138451
$ws.Range("O8").Value = 0.04138006727350978
$ws.Range("P8").Value = 0.04138006727350978
$ws.Range("Q8").Value = 1.310125159215
$ws.Range("R8").Value = 11.791126432935
$ws.Range("S8").Value = 0.008379832468295398
$ws.Range("T8").Value = 0.008379832468295394
$ws.Range("G9").Value = 0.7648950000000001
$ws.Range("H9").Value = 2.294685
$ws.Range("I9").Value = 0.2025089136010154
$ws.Range("J9").Value = 0.2025089136010154
$ws.Range("M9").Value = 0.02775933333333333
$ws.Range("N9").Value = 0.083278
$ws.Range("O9").Value = 0.000670639701031176
$ws.Range("P9").Value = 0.000670639701031176
$ws.Range("Q9").Value = 0.02123297527
$ws.Range("R9").Value = 0.19109677743
$ws.Range("S9").Value = 0.0001358105172735332
$ws.Range("T9").Value = 0.0001358105172735332
